# Auto-generated Excel COM-interop script to apply value updates
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each change corresponds to a cell value update in the source OOXML diff;
# some cells are newly added (previously empty) or fully cleared (removed).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 27950
$ws.Range("J3").Value = 27950
$ws.Range("L3").Value = 27950
$ws.Range("N3").Value = -28178
$ws.Range("H9").Value = 231.6
$ws.Range("I9").Value = 303.22223
$ws.Range("K9").Value = 303.22223
$ws.Range("M9").Value = -134.22223
$ws.Range("H40").Value = 4103.159
$ws.Range("I40").Value = 3071.4324
$ws.Range("K40").Value = 3071.4324
$ws.Range("M40").Value = -2896.4324
$ws.Range("H102").Value = 27950
$ws.Range("J102").Value = 27950
$ws.Range("L102").Value = 27950
$ws.Range("N102").Value = -34440
$ws.Range("H113").Value = 4998.3335
$ws.Range("J113").Value = 6666.6665
$ws.Range("L113").Value = 6666.6665
$ws.Range("N113").Value = -13174.6665
$ws.Range("H132").Value = 15393.429
$ws.Range("I132").Value = 15393.429
$ws.Range("K132").Value = 46180.287
$ws.Range("M132").Value = -43650.287
$ws.Range("H137").Value = 2461.7036
$ws.Range("I137").Value = 1516.1333
$ws.Range("K137").Value = 4548.3999
$ws.Range("M137").Value = -1998.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2038.72
$ws.Range("I32").Value = 900.6875
$ws.Range("K32").Value = 900.6875
$ws.Range("M32").Value = -613.6875
$ws.Range("H61").Value = 3487.889
$ws.Range("I61").Value = 3487.889
$ws.Range("K61").Value = 3487.889
$ws.Range("M61").Value = -3275.889
$ws.Range("H74").Value = 4760.5
$ws.Range("I74").Value = 3941.5715
$ws.Range("J74").Value = 6671.3335
$ws.Range("K74").Value = 3941.5715
$ws.Range("L74").Value = 6671.3335
$ws.Range("M74").Value = -3067.5715
$ws.Range("N74").Value = -8419.333500000001
$ws.Range("H77").Value = 4760.5
$ws.Range("I77").Value = 3941.5715
$ws.Range("J77").Value = 6671.3335
$ws.Range("K77").Value = 19707.8575
$ws.Range("L77").Value = 33356.6675
$ws.Range("M77").Value = -15339.8575
$ws.Range("N77").Value = -42092.6675
$ws.Range("H110").Value = 845.2
$ws.Range("I110").Value = 828
$ws.Range("J110").Value = 885.3333
$ws.Range("K110").Value = 828
$ws.Range("L110").Value = 885.3333
$ws.Range("M110").Value = 1217
$ws.Range("N110").Value = -4975.3333
$ws.Range("H132").Value = 5474.0835
$ws.Range("I132").Value = 5244.4546
$ws.Range("K132").Value = 15733.3638
$ws.Range("M132").Value = -13203.3638
$ws.Range("H136").Value = 3487.889
$ws.Range("I136").Value = 3487.889
$ws.Range("K136").Value = 10463.667
$ws.Range("M136").Value = -7913.667000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3341.6667
$ws.Range("H22").Value = 1701.375
$ws.Range("I22").Value = 1701.375
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1701.375
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1528.375
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H134").Value = 2500
$ws.Range("I134").Value = 2500
$ws.Range("K134").Value = 7500
$ws.Range("M134").Value = -4965
$ws.Range("N22").ClearContents()
$ws.Range("M64").ClearContents()
$ws.Range("M67").ClearContents()
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5815.391
$ws.Range("I31").Value = 2880.6667
$ws.Range("K31").Value = 2880.6667
$ws.Range("M31").Value = -2585.6667
$ws.Range("H34").Value = 5815.391
$ws.Range("I34").Value = 2880.6667
$ws.Range("K34").Value = 2880.6667
$ws.Range("M34").Value = -2678.6667
$ws.Range("H58").Value = 2985.0625
$ws.Range("I58").Value = 2394.8333
$ws.Range("K58").Value = 2394.8333
$ws.Range("M58").Value = -2191.8333
$ws.Range("H107").Value = 335.58334
$ws.Range("I107").Value = 173
$ws.Range("J107").Value = 563.2
$ws.Range("K107").Value = 173
$ws.Range("L107").Value = 563.2
$ws.Range("M107").Value = 1747
$ws.Range("N107").Value = -4403.2
$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470
$ws.Range("H136").Value = 2985.0625
$ws.Range("I136").Value = 2394.8333
$ws.Range("K136").Value = 7184.499899999999
$ws.Range("M136").Value = -4634.499899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1824.9
$ws.Range("J68").Value = 1806.125
$ws.Range("L68").Value = 5418.375
$ws.Range("N68").Value = -7040.375
$ws.Range("H71").Value = 1824.9
$ws.Range("J71").Value = 1806.125
$ws.Range("L71").Value = 16255.125
$ws.Range("N71").Value = -24367.125
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("H107").Value = 288.25
$ws.Range("I107").Value = 301
$ws.Range("K107").Value = 903
$ws.Range("M107").Value = 1017
$ws.Range("H131").Value = 1718
$ws.Range("I131").Value = 1017.55554
$ws.Range("K131").Value = 3052.66662
$ws.Range("M131").Value = 1987.33338
$ws.Range("H140").Value = 3208.3333
$ws.Range("I140").Value = 2850
$ws.Range("K140").Value = 8550
$ws.Range("M140").Value = -3370
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("H63").Value = 50000
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("H66").Value = 50000
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H132").Value = 8964.666999999999
$ws.Range("I132").Value = 8947.5
$ws.Range("K132").Value = 26842.5
$ws.Range("M132").Value = -24312.5
$ws.Range("N59").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("M63").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("M66").ClearContents()
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1076
$ws.Range("I16").Value = 1076
$ws.Range("K16").Value = 1076
$ws.Range("M16").Value = -906
$ws.Range("H46").Value = 5954.1875
$ws.Range("I46").Value = 4255.6665
$ws.Range("J46").Value = 8138
$ws.Range("K46").Value = 4255.6665
$ws.Range("L46").Value = 8138
$ws.Range("M46").Value = -4067.6665
$ws.Range("N46").Value = -8514
$ws.Range("H63").Value = 44444
$ws.Range("I63").Value = 44444
$ws.Range("K63").Value = 44444
$ws.Range("M63").Value = -43695
$ws.Range("H66").Value = 44444
$ws.Range("I66").Value = 44444
$ws.Range("K66").Value = 133332
$ws.Range("M66").Value = -129588
$ws.Range("H68").Value = 6890.636
$ws.Range("I68").Value = 3474.25
$ws.Range("J68").Value = 8842.857
$ws.Range("K68").Value = 3474.25
$ws.Range("L68").Value = 8842.857
$ws.Range("M68").Value = -2725.25
$ws.Range("N68").Value = -10340.857
$ws.Range("H71").Value = 6890.636
$ws.Range("I71").Value = 3474.25
$ws.Range("J71").Value = 8842.857
$ws.Range("K71").Value = 17371.25
$ws.Range("L71").Value = 44214.285
$ws.Range("M71").Value = -13627.25
$ws.Range("N71").Value = -51702.285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 22568.25
$ws.Range("J41").Value = 21761.25
$ws.Range("L41").Value = 21761.25
$ws.Range("N41").Value = -22541.25
$ws.Range("H136").Value = 2830.3333
$ws.Range("I136").Value = 2337.2273
$ws.Range("K136").Value = 7011.6819
$ws.Range("M136").Value = -4461.6819
